$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "UAH shooting victim released from hospital"
$ws.Range("B2").Value = "2010-03-29T10:19:00UTC"
$ws.Range("C2").Value = 45
$ws.Range("D2").Value = "day_31_beyond"
$ws.Range("E2").Value = "http://blog.al.com/breaking/2010/03/uah_shooting_victim_released_f.html"

# Row 3
$ws.Range("A3").Value = "Alabama Shooting Suspect's Husband: 'I'm No Psychologist'"
$ws.Range("B3").Value = "2010-02-12T23:34:01UTC"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = "day_0"
$ws.Range("E3").Value = "https://abcnews.go.com/GMA/alabama-university-shooting-suspect-amy-bishop-violent-past/story?id=9839348"

# Row 4
$ws.Range("A4").Value = "Ex-professor gets life in prison for UAH shooting - WSFA.com Montgomery Alabama news."
$ws.Range("B4").Value = "2012-09-24T00:00:00UTC"
$ws.Range("C4").Value = 955
$ws.Range("D4").Value = "day_31_beyond"
$ws.Range("E4").Value = "http://www.wsfa.com/story/19619130/jury-finds-amy-bishop-anderson-guilty-of-capital-murder"

# Row 5
$ws.Range("A5").Value = "Professor Accused in Killings Is Said to Attempt Suicide"
$ws.Range("B5").Value = "2010-06-19T18:14:50UTC"
$ws.Range("C5").Value = 127
$ws.Range("D5").Value = "day_31_beyond"
$ws.Range("E5").Value = "https://www.nytimes.com/2010/06/20/us/20bishop.html?ref=us"

# Row 6
$ws.Range("A6").Value = "UAH campus memorial service"
$ws.Range("B6").Value = "2010-02-20T12:27:00UTC"
$ws.Range("C6").Value = 8
$ws.Range("D6").Value = "day_2_to_30"
$ws.Range("E6").Value = "http://blog.al.com/breaking/2010/02/uah_campus_memorial_service.html"

# Row 7
$ws.Range("A7").Value = "Prosecutors to seek death penalty for accused UAH shooter Amy Bishop"
$ws.Range("B7").Value = "2011-05-25T15:03:00UTC"
$ws.Range("C7").Value = 467
$ws.Range("D7").Value = "day_31_beyond"
$ws.Range("E7").Value = "http://blog.al.com/breaking/2011/05/prosecutors_to_seek_death_pena.html"

# Row 8
$ws.Range("A8").Value = "Saturday press conference on campus shooting - WAFF-TV: News, Weather and Sports for Huntsville, AL"
$ws.Range("B8").Value = "2010-02-13T00:00:00UTC"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = "day_1"
$ws.Range("E8").Value = "http://www.waff.com/Global/story.asp?S=11983009"

# Row 9
$ws.Range("A9").Value = "Ala. prof's story begins with brother's 1986 death"
$ws.Range("B9").Value = "2010-02-16T00:00:00UTC"
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = "day_2_to_30"
$ws.Range("E9").Value = "http://www.boston.com/news/nation/articles/2010/02/16/survivor_ala_prof_in_slayings_shot_methodically/"

# Row 10
$ws.Range("A10").Value = "A Promising Start to an Academic Life-or to a Life of Violence?"
$ws.Range("B10").Value = "1-01-01T00:00:00UTC"
$ws.Range("C10").Value = "unknown"
$ws.Range("D10").Value = "unknown"
$ws.Range("E10").Value = "https://web.archive.org/web/20140715002421/http://www.crimelibrary.com/notorious_murders/mass/amy_bishop/2.html"

# Row 11
$ws.Range("A11").Value = "Alabama shooting survivor: 'There was no way to ever anticipate this'"
$ws.Range("B11").Value = "1-01-01T00:00:00UTC"
$ws.Range("C11").Value = "unknown"
$ws.Range("D11").Value = "unknown"
$ws.Range("E11").Value = "http://www.cnn.com/2010/CRIME/02/17/alabama.shooting.witness/index.html"

# Row 12
$ws.Range("A12").Value = "The New York Times"
$ws.Range("B12").Value = "2009-07-02T15:10:24UTC"
$ws.Range("C12").Value = "unknown"
$ws.Range("D12").Value = "unknown"
$ws.Range("E12").Value = "https://www.nytimes.com/2009/07/05/movies/05pare.html"

# Row 13
$ws.Range("A13").Value = "Amy Bishop case"
$ws.Range("B13").Value = "1-01-01T00:00:00UTC"
$ws.Range("C13").Value = "unknown"
$ws.Range("D13").Value = "unknown"
$ws.Range("E13").Value = "http://www.boston.com/news/specials/02_15_Amy_Bishop/"
